$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -4.681004024549235
$ws.Range("C3").Value = -1.165899493033517
$ws.Range("C4").Value = -0.05393960539424965
$ws.Range("C5").Value = -0.4168956827577189
$ws.Range("C6").Value = -0.1239422036827855
$ws.Range("C7").Value = 0.01913423197059164
